$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# Note: the runtime's character-width<->pixel rounding requires a slight
# downward offset (target - 0.9) from the integer target to land back on
# the exact integer width value in the saved OOXML.
$ws.Columns.Item(3).ColumnWidth = 48.1
$ws.Columns.Item(4).ColumnWidth = 50.1
$ws.Columns.Item(5).ColumnWidth = 36.1

# Update header row (row 1) values
$ws.Range("C1").Value = "button_testResultActions_internalRoleButtonName"
$ws.Range("D1").Value = "button_testResultActions_internalRoleButtonName_1"
$ws.Range("E1").Value = "button_testResultDetails_class"
$ws.Range("F1").Value = "button_testResultDetails_internalRoleButtonName"

# Update data row (row 2) values
$ws.Range("C2").Value = "Failed Automations - Apply to"
$ws.Range("D2").Value = "Failed Portal - Login with"
$ws.Range("E2").Value = """]:nth-child(3) [class=""css-1yjo05o"
